$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume data rows (values refreshed by the scraper)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '24.506.69'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.38%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.669.27'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.39%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9987'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.45%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.18'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.88%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9994'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.25%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3969'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.45%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3920'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.37%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '51.98'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +6.85%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.410'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.47%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9986'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.45%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08614'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.54%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '24.47'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.18%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.357'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.49%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001352'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.84%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.922'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +5.06%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.668.04'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.32%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '95.54'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.15%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06984'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.56%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '20.69'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.89%  '

$ws.Range("E21").Value = '  +0.76%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9986'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.30%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.78'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.29%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '24.490.34'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.24%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.428'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.67%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.058'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +11.60%  '

$ws.Range("E27").Value = '  -0.03%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '157.66'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.41%  '

$ws.Range("B29").Value = 'BitcoinCash'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '143.07'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.69%  '

$ws.Range("B30").Value = 'HuobiToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.471'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.49%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.178'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -10.29%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.552'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.82%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.849.29'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.23%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.064'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +7.49%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.08303'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.66%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.03039'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.64%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.910'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -5.11%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '11.32'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +12.26%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2780'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.28%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.09261'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.34%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.7780'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.57%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '13.93'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.99%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.448'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.06%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.62'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.96%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.7154'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.38%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.547'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.84%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.150'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.35%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.9991'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.23%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.08463'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.29%  '

$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '136.84'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.79%  '

$ws.Range("B51").Value = 'Flow'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QQ0NCmjVq+flow-flow'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.280'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.65%  '
